# Regenerate orders with updated distance/sizes.
# Global rename (applied to the text content of every cell, wherever the
# token appears as a substring):
#   D51 -> D55
#   D64 -> D69
#   D80 -> D86
#   S30 -> S31
# (S20 / S25 are left untouched.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ur = $ws.UsedRange
$startRow = $ur.Row()
$startCol = $ur.Column()
$endRow = $startRow + $ur.Rows.Count - 1
$endCol = $startCol + $ur.Columns.Count - 1

for ($r = $startRow; $r -le $endRow; $r++) {
    for ($c = $startCol; $c -le $endCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()

        if ($val -is [string]) {
            if ($val.Contains("D51") -or $val.Contains("D64") -or $val.Contains("D80") -or $val.Contains("S30")) {
                $newVal = $val.Replace("D51", "D55").Replace("D64", "D69").Replace("D80", "D86").Replace("S30", "S31")
                $cell.Value = $newVal
            }
        }
    }
}
